$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = "2024-10-05 00:00:00"
$ws.Range("B66").Value = 75650
$ws.Range("C66").Value = 10756.89
$ws.Range("D66").Value = 9519.370000000001
$ws.Range("E66").Value = 7.0184
